$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------------
# 1. Copy row formatting (styles) from existing template rows onto the new
#    rows so that the cellXf (style) indices match exactly what Excel would
#    naturally produce for these kinds of rows.
# ---------------------------------------------------------------------------

# Section header rows (THURSDAY / FRIDAY) - template: row 5 (THURSDAY header)
$ws.Range("A5:F5").Copy()
$ws.Range("A562").PasteSpecial(-4122)
$ws.Range("A568").PasteSpecial(-4122)

# Lockup rows (ht=30) - template: row 555 (Lockup row)
$ws.Range("A555:F555").Copy()
$ws.Range("A563").PasteSpecial(-4122)
$ws.Range("A564").PasteSpecial(-4122)

# AV Shutdown rows with ht=45 - template: row 6
$ws.Range("A6:F6").Copy()
$ws.Range("A569").PasteSpecial(-4122)
$ws.Range("A570").PasteSpecial(-4122)

# AV Shutdown rows with no explicit row height - template: row 548
$ws.Range("A548:F548").Copy()
$ws.Range("A571").PasteSpecial(-4122)
$ws.Range("A572").PasteSpecial(-4122)
$ws.Range("A573").PasteSpecial(-4122)
$ws.Range("A574").PasteSpecial(-4122)
$ws.Range("A575").PasteSpecial(-4122)
$ws.Range("A576").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Set row heights to match the source rows (Excel auto-fits these normally
#    because of wrapped text, but the headless engine needs an explicit nudge)
# ---------------------------------------------------------------------------
$ws.Rows.Item(563).RowHeight = 30
$ws.Rows.Item(564).RowHeight = 30
$ws.Rows.Item(569).RowHeight = 45
$ws.Rows.Item(570).RowHeight = 45

# ---------------------------------------------------------------------------
# 3. Fill in the cell values for the new rows.
# ---------------------------------------------------------------------------

# Row 562 - THURSDAY section header
$ws.Range("B562").Value = "THURSDAY"

# Row 563 - Lockup, Thu Oct 27 2016
$ws.Range("A563").Value = "Lockup"
$ws.Range("B563").Value = 42670
$ws.Range("C563").Value = "1900"
$ws.Range("D563").Value = "CLH"
$ws.Range("E563").Value = "J"
$ws.Range("F563").Value = "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS."

# Row 564 - Lockup, Thu Oct 27 2016
$ws.Range("A564").Value = "Lockup"
$ws.Range("B564").Value = 42670
$ws.Range("C564").Value = "1900"
$ws.Range("D564").Value = "CLH"
$ws.Range("E564").Value = "M"
$ws.Range("F564").Value = "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS."

# Row 568 - FRIDAY section header
$ws.Range("B568").Value = "FRIDAY"

# Row 569 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A569").Value = "AV Shutdown"
$ws.Range("B569").Value = 42671
$ws.Range("C569").Value = "1600"
$ws.Range("D569").Value = "BC"
$ws.Range("E569").Value = "320"
$ws.Range("F569").Value = "Turn off PC and Projector in room. Projector remote is on PC cart - please leave on PC cart. ALL EQUIPMENT STAYS IN ROOM. Lock room - key on Bethune classroom keys in CB 121A storeroom."

# Row 570 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A570").Value = "AV Shutdown"
$ws.Range("B570").Value = 42671
$ws.Range("C570").Value = "1600"
$ws.Range("D570").Value = "LUM"
$ws.Range("E570").Value = "306"
$ws.Range("F570").Value = "Turn off PC and Projector in room. ALL EQUIPMENT STAYS IN ROOM. Lock room. Get key from Dean's office in 3rd floor Lumbers."

# Row 571 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A571").Value = "AV Shutdown"
$ws.Range("B571").Value = 42671
$ws.Range("C571").Value = "1700"
$ws.Range("D571").Value = "CB"
$ws.Range("E571").Value = "121"
$ws.Range("F571").Value = "Return neck mic to drawer."

# Row 572 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A572").Value = "AV Shutdown"
$ws.Range("B572").Value = 42671
$ws.Range("C572").Value = "1700"
$ws.Range("D572").Value = "LSB"
$ws.Range("E572").Value = "101"
$ws.Range("F572").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

# Row 573 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A573").Value = "AV Shutdown"
$ws.Range("B573").Value = 42671
$ws.Range("C573").Value = "1700"
$ws.Range("D573").Value = "LSB"
$ws.Range("E573").Value = "103"
$ws.Range("F573").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

# Row 574 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A574").Value = "AV Shutdown"
$ws.Range("B574").Value = 42671
$ws.Range("C574").Value = "1700"
$ws.Range("D574").Value = "LSB"
$ws.Range("E574").Value = "105"
$ws.Range("F574").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

# Row 575 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A575").Value = "AV Shutdown"
$ws.Range("B575").Value = 42671
$ws.Range("C575").Value = "1700"
$ws.Range("D575").Value = "LSB"
$ws.Range("E575").Value = "106"
$ws.Range("F575").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

# Row 576 - AV Shutdown, Fri Oct 28 2016
$ws.Range("A576").Value = "AV Shutdown"
$ws.Range("B576").Value = 42671
$ws.Range("C576").Value = "1700"
$ws.Range("D576").Value = "LSB"
$ws.Range("E576").Value = "107"
$ws.Range("F576").Value = "Make sure neck mic goes back to drawer and log off touchscreen."

# ---------------------------------------------------------------------------
# 4. Update the view state to match: scrolled so row 558 is at the top, and
#    F581 is the active/selected cell.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 558
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F581").Select()
